$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New retailer name added in row 62 (B column), reusing existing "WALIDAD" text
$ws.Range("B62").Value = "WALIDAD"

# New column Q values added for several retailers (Dec 12 daily update)
$ws.Range("Q3").Value = 5200
$ws.Range("Q4").Value = 3120
$ws.Range("Q5").Value = 1040
$ws.Range("Q6").Value = 2080
$ws.Range("Q7").Value = 2080
$ws.Range("Q25").Value = 3120
$ws.Range("Q31").Value = 5200
$ws.Range("Q57").Value = 2080
$ws.Range("Q61").Value = 2080
$ws.Range("Q65").Value = 5200
$ws.Range("Q78").Value = 5200
$ws.Range("Q80").Value = 5200
$ws.Range("Q82").Value = 2080
$ws.Range("Q91").Value = 1040
